$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "30.034.24"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").Value = "2.103.09"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.65%  "
$ws.Range("D5").Value = "348.26"
$ws.Range("E5").Value = "  +3.35%  "
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("D7").Value = "0.5159"
$ws.Range("E7").Value = "  -1.65%  "
$ws.Range("D8").Value = "0.4440"
$ws.Range("E8").Value = "  -2.63%  "
$ws.Range("D9").Value = "52.36"
$ws.Range("E9").Value = "  -4.30%  "
$ws.Range("D10").Value = "0.08969"
$ws.Range("D11").Value = "1.170"
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").Value = "25.51"
$ws.Range("E12").Value = "  +3.97%  "
$ws.Range("D13").Value = "2.107.67"
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("D14").Value = "8.252"
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("D15").Value = "6.738"
$ws.Range("E15").Value = "  -1.78%  "
$ws.Range("D16").Value = "99.13"
$ws.Range("E16").Value = "  +2.09%  "
$ws.Range("D17").Value = "0.00001146"
$ws.Range("E17").Value = "  -2.79%  "
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("D19").Value = "21.03"
$ws.Range("E19").Value = "  +8.15%  "
$ws.Range("D20").Value = "0.06678"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("D22").Value = "6.235"
$ws.Range("E22").Value = "  -1.07%  "
$ws.Range("D23").Value = "30.143.56"
$ws.Range("E23").Value = "  -1.67%  "
$ws.Range("D24").Value = "12.69"
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("D25").Value = "2.336"
$ws.Range("E25").Value = "  -1.02%  "
$ws.Range("D26").Value = "2.354.20"
$ws.Range("E26").Value = "  -0.70%  "
$ws.Range("D27").Value = "21.96"
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").Value = "162.30"
$ws.Range("E29").Value = "  -1.30%  "
$ws.Range("D30").Value = "133.55"
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("D31").Value = "1.175"
$ws.Range("E31").Value = "  -3.08%  "
$ws.Range("D32").Value = "0.1064"
$ws.Range("E32").Value = "  -0.95%  "
$ws.Range("D33").Value = "1.640"
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("D34").Value = "6.218"
$ws.Range("E34").Value = "  -2.45%  "
$ws.Range("D35").Value = "3.970"
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("D36").Value = "5.912"
$ws.Range("E36").Value = "  +0.47%  "
$ws.Range("D37").Value = "10.21"
$ws.Range("E37").Value = "  -3.92%  "
$ws.Range("D38").Value = "0.02574"
$ws.Range("E38").Value = "  -2.40%  "
$ws.Range("D39").Value = "0.06790"
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("D40").Value = "0.2295"
$ws.Range("E40").Value = "  -1.49%  "
$ws.Range("D41").Value = "12.52"
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("D42").Value = "0.6800"
$ws.Range("E42").Value = "  -1.30%  "
$ws.Range("D43").Value = "1.295"
$ws.Range("E43").Value = "  +2.73%  "
$ws.Range("D44").Value = "14.30"
$ws.Range("E44").Value = "  -2.78%  "
$ws.Range("D45").Value = "0.6362"
$ws.Range("E45").Value = "  -1.95%  "
$ws.Range("D46").Value = "2.287"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("D47").Value = "0.00000000360"
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("D48").Value = "3.634"
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("D49").Value = "1.219"
$ws.Range("E49").Value = "  -2.95%  "
$ws.Range("D50").Value = "82.39"
$ws.Range("E50").Value = "  -1.35%  "
$ws.Range("E51").Value = "  +0.37%  "
